$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": refresh the top-level P&L / trade-count metrics now that
# trade #76 has closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.13   # Current Capital
$summary.Range("B4").Value = 0.12      # Total P&L $
$summary.Range("B5").Value = 0.03      # Total P&L %
$summary.Range("B6").Value = 76        # Total Trades
$summary.Range("B8").Value = 24        # Losing Trades
$summary.Range("B9").Value = 43.42     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": mirror the same refresh for the MarketMaking row.
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.13     # Capital
$status.Range("D4").Value = 76         # Trades
$status.Range("E4").Value = 0.12       # P&L $
$status.Range("F4").Value = 0.13       # P&L %
$status.Range("G4").Value = 43.42      # Win Rate %

# ---------------------------------------------------------------------------
# Append the newly-closed trade #76 as row 77 on both the "All Trades" and
# "MarketMaking" logs (they track the same data in this workbook).
# ---------------------------------------------------------------------------
$tradeSheets = @("All Trades", "MarketMaking")
foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A77").Value = 76

    # Dates in this log are stored as plain text (e.g. "2026-02-17"), not as
    # real Excel date serials, so force the cell to Text before assigning
    # the value to stop Excel's autodetection from converting it.
    $ws.Range("B77").NumberFormat = "@"
    $ws.Range("B77").Value = "2026-02-17"
    $ws.Range("B77").Style = "Normal"

    $ws.Range("C77").Value = "12:57:39"
    $ws.Range("D77").Value = "MarketMaking"
    $ws.Range("E77").Value = "UP"
    $ws.Range("F77").Value = 0.58
    $ws.Range("G77").Value = 0.48
    $ws.Range("H77").Value = "CLOSED"
    $ws.Range("I77").Value = -17.2414
    $ws.Range("J77").Value = -0.1
    $ws.Range("K77").Value = 100.13
    $ws.Range("L77").Value = 0
    $ws.Range("M77").Value = 0
    $ws.Range("N77").Value = 0.6
    $ws.Range("O77").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P77").Value = "early_exit"
    $ws.Range("Q77").Value = 0.13
}
